# Append 5 new applicant rows (141-145) to the "Лист1" sheet, matching the
# source data added on 2025-07-01 (commit: "18:11 time 01.07.2025 date").
#
# Columns: A=F.I.Sh, B=Yo'nalish, C=Ta'lim tili, D=Ta'lim shakli,
#          E=Passport, F=JSHIR, G=Viloyat, H=Tuman,
#          I=Telegram raqami, J=Telefon raqami, K=Sana
#
# Columns F, I, J, K hold values that look numeric/date-like
# (ID numbers, phone numbers, ISO dates). In the source workbook these are
# stored as plain text (inlineStr) cells with no special number format, so
# we must force Excel to keep them as text instead of silently converting
# them to numbers / date serials. We do this the same way a user typing
# into Excel would (leading apostrophe => quote-prefixed text), and then
# immediately re-apply the plain formatting of a known "clean" text cell
# (the previous row) on top, so no stray number-format/style gets baked in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=141; A="Azamov Aziz"; B="Hayot faoliyati xavfsizligi"; C="O'zbek tili"; D="Kunduzgi"; E="AD5352798"; F="52810076230020"; G="Surxondaryo viloyati"; H="Denov tumani"; I="998938239646"; J="+998938239646"; K="2025-07-01" },
    @{ Row=142; A="XASANOVA SEVINCHXON AXRORBEK QIZI"; B="Yurisprudensiya"; C="Rus tili"; D="Kunduzgi"; E="AC2540058"; F="60810036940010"; G="Toshkent shahri"; H="Yangihayot tumani"; I="998944941008"; J="+998944941008"; K="2025-07-01" },
    @{ Row=143; A="Sayfiddinova Maqsad Sherzodovna"; B="Yurisprudensiya"; C="O'zbek tili"; D="Kunduzgi"; E="AD5715114"; F="60910076150056"; G="Samarqand viloyati"; H="Toyloq tumani"; I="998948370910"; J="+998993507867"; K="2025-07-01" },
    @{ Row=144; A="Bekmurodov Sardor Shuhratovich"; B="Yurisprudensiya"; C="O'zbek tili"; D="Kunduzgi"; E="AD9500428"; F="51608076270020"; G="Surxondaryo viloyati"; H="Termiz shahri"; I="998996212007"; J="+998942003007"; K="2025-07-01" },
    @{ Row=145; A="Po'latov Fayzulloh Kamoliddin o'g'li"; B="Menejment"; C="O'zbek tili"; D="Kunduzgi"; E="AD8066834"; F="51607086540127"; G="Toshkent shahri"; H="Shayxontohur tumani"; I="998994778754"; J="+998949928754"; K="2025-07-01" }
)

$xlPasteFormats = -4122

# Row 140 is the last existing data row and has plain (unstyled) cells for
# every column -- use it as the initial "clean format" donor.
$donorRow = 140

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Formula = "'" + $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Formula = "'" + $r.I
    $ws.Cells.Item($row, 10).Formula = "'" + $r.J
    $ws.Cells.Item($row, 11).Formula = "'" + $r.K

    # Re-apply plain formatting (no quote-prefix / number-format side
    # effects) from the donor row, column by column, without touching the
    # values we just set.
    $ws.Range("F$donorRow").Copy()
    $ws.Range("F$row").PasteSpecial($xlPasteFormats)

    $ws.Range("I$donorRow" + ":J$donorRow").Copy()
    $ws.Range("I$row" + ":J$row").PasteSpecial($xlPasteFormats)

    $ws.Range("K$donorRow").Copy()
    $ws.Range("K$row").PasteSpecial($xlPasteFormats)

    $donorRow = $row
}

$excel.CutCopyMode = 0

Write-Host "Added rows 141-145"
